{"js": "// Module 3.2 paper edits:\n//  1) The paragraph that used to read\n//       \"Post-2021 resources, such as blogs and articles from industry\n//        leaders, highlight trends influenced by DevOps and agile\n//        methodologies:\"\n//     now reads\n//       \"Some of the highlight trends influenced by DevOps and agile\n//        methodologies are:\"\n//  2) Inside the \"Comparison and Contrast\" paragraph, \"recent trends\" was\n//     changed to \"current trends\".\n\nconst body = context.document.body;\n\n// --- Edit 1 -------------------------------------------------------------\nconst introResults = body.search(\n  \"Post-2021 resources, such as blogs and articles from industry leaders, highlight trends influenced by DevOps and agile methodologies:\",\n  { matchCase: true }\n);\nintroResults.load(\"items\");\nawait context.sync();\n\nif (introResults.items.length > 0) {\n  introResults.items[0].insertText(\n    \"Some of the highlight trends influenced by DevOps and agile methodologies are:\",\n    \"Replace\"\n  );\n}\n\n// --- Edit 2 -------------------------------------------------------------\nconst trendsResults = body.search(\"recent trends\", { matchCase: true });\ntrendsResults.load(\"items\");\nawait context.sync();\n\nif (trendsResults.items.length > 0) {\n  trendsResults.items[0].insertText(\"current trends\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Module 3.2 paper edits:\n#  1) The paragraph that used to read\n#       \"Post-2021 resources, such as blogs and articles from industry\n#        leaders, highlight trends influenced by DevOps and agile\n#        methodologies:\"\n#     now reads\n#       \"Some of the highlight trends influenced by DevOps and agile\n#        methodologies are:\"\n#  2) Inside the \"Comparison and Contrast\" paragraph, \"recent trends\" was\n#     changed to \"current trends\".\n\n$d = $word.ActiveDocument\n\n# --- Edit 1 ---------------------------------------------------------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$findText1 = \"Post-2021 resources, such as blogs and articles from industry leaders, highlight trends influenced by DevOps and agile methodologies:\"\n$replaceText1 = \"Some of the highlight trends influenced by DevOps and agile methodologies are:\"\n$find1.Execute($findText1, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText1, 2)\n\n# --- Edit 2 -----------------------------------------------------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"recent trends\", $false, $false, $false, $false, $false, $true, 1, $false, \"current trends\", 2)\n"}
